# This workbook's single data sheet lists weekly Fruta/Hortaliza price
# observations (Macroferia Regional de Talca - Granada). The update re-syncs
# the 24 data rows (rows 2-25) against the latest source extract: each row's
# date/quality/volume/price/unit/origin fields get replaced by the values
# that another (pre-existing) row used to hold - i.e. the 24 rows are
# permuted among themselves across columns D and L:T, while columns
# A:C, E:K stay identical (they are constant for this whole sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row -> source row (values currently sitting in $source row move to $newRow)
$mapping = @{
    2  = 9
    3  = 16
    4  = 3
    5  = 19
    6  = 20
    7  = 22
    8  = 2
    9  = 11
    10 = 7
    11 = 17
    12 = 6
    13 = 21
    14 = 18
    15 = 10
    16 = 4
    17 = 12
    18 = 23
    19 = 8
    20 = 14
    21 = 13
    22 = 24
    23 = 25
    24 = 5
    25 = 15
}

# Columns that actually vary row-to-row and therefore need to move.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot every source cell's value BEFORE any writes, so overlapping
# writes during the permutation don't clobber data we still need to read.
$snapshot = @{}
for ($r = 2; $r -le 25; $r++) {
    foreach ($col in $cols) {
        $snapshot["$col$r"] = $ws.Range("$col$r").Value2
    }
}

for ($newRow = 2; $newRow -le 25; $newRow++) {
    $srcRow = $mapping[$newRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $snapshot["$col$srcRow"]
    }
}
